# Apply the "departamentos / areas comunes" relabeling to the condominios
# report template:
#   - J6/K6 header labels change ("Total Usuarios" / "N.\xB0 Torres" ->
#     "Total Departamentos" / "Areas Comunes"); the old L6 header
#     ("N.\xB0 Bloques de Casas") is dropped entirely.
#   - Row 7 (the template/placeholder row) shifts left to match: the last
#     placeholder (cantidadCasas) is dropped and the remaining two are
#     renamed to totalDepartamentos / areasComunes.
#   - Column widths for J/K grow to fit the new labels.
#   - Selection moves to K6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 6) -------------------------------------------------
$ws.Range("J6").Value = "Total Departamentos"
$ws.Range("K6").Value = "Areas Comunes"

# L6 used to hold "N.\xB0 Bloques de Casas"; it is now empty with the
# worksheet's plain default formatting (no fill/border/bold like the other
# header cells). Clear it, then pull default formatting from a pristine,
# never-touched cell so L6 truly resets instead of keeping the old header
# style.
$ws.Range("L6").ClearContents()
$ws.Range("ZZ6").Copy()
$ws.Range("L6").PasteSpecial(-4122)

# --- Template row (row 7) ------------------------------------------------
$ws.Range("J7").Value = "`${table:condominios.totalDepartamentos}"
$ws.Range("K7").Value = "`${table:condominios.areasComunes}"

$ws.Range("L7").ClearContents()
$ws.Range("ZZ7").Copy()
$ws.Range("L7").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Column widths ---------------------------------------------------
$ws.Columns.Item(10).ColumnWidth = 22.95
$ws.Columns.Item(11).ColumnWidth = 17.67

# --- Selection / scroll position ---------------------------------------
$ws.Range("K6").Select()
$excel.ActiveWindow.ScrollColumn = 3
